$wb = $excel.ActiveWorkbook

# --- ALC sheet updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 90.818184
$ws.Range("I12").Value = 104
$ws.Range("J12").Value = 75
$ws.Range("K12").Value = 104
$ws.Range("L12").Value = 75
$ws.Range("M12").Value = 66
$ws.Range("N12").Value = -415
$ws.Range("H93").Value = 44250
$ws.Range("J93").Value = 44250
$ws.Range("L93").Value = 44250
$ws.Range("N93").Value = -49242
$ws.Range("H137").Value = 2083.5715
$ws.Range("I137").Value = 2044.4736
$ws.Range("J137").Value = 2166.111
$ws.Range("K137").Value = 6133.4208
$ws.Range("L137").Value = 6498.333
$ws.Range("M137").Value = -3583.4208
$ws.Range("N137").Value = -11598.333
$ws.Range("H138").Value = 6582027.5
$ws.Range("I138").Value = 1592.6666
$ws.Range("J138").Value = 7815859
$ws.Range("K138").Value = 4777.9998
$ws.Range("L138").Value = 23447577
$ws.Range("M138").Value = 362.0002000000004
$ws.Range("N138").Value = -23457857
$ws.Range("H141").Value = 2939
$ws.Range("I141").Value = 2886.25
$ws.Range("J141").Value = 3150
$ws.Range("K141").Value = 8658.75
$ws.Range("L141").Value = 9450
$ws.Range("M141").Value = -3478.75
$ws.Range("N141").Value = -19810

# --- ARM sheet updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 779.8
$ws.Range("I4").Value = 300
$ws.Range("J4").Value = 1099.6666
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 1099.6666
$ws.Range("M4").Value = -184
$ws.Range("N4").Value = -1331.6666
$ws.Range("H32").Value = 7715.6455
$ws.Range("I32").Value = 3709.7014
$ws.Range("J32").Value = 30082.166
$ws.Range("K32").Value = 3709.7014
$ws.Range("L32").Value = 30082.166
$ws.Range("M32").Value = -3422.7014
$ws.Range("N32").Value = -30656.166
$ws.Range("H45").Value = 6130.4287
$ws.Range("I45").Value = 7524.5625
$ws.Range("K45").Value = 7524.5625
$ws.Range("M45").Value = -7147.5625
$ws.Range("H74").Value = 50041.332
$ws.Range("I74").Value = 64127.75
$ws.Range("K74").Value = 64127.75
$ws.Range("M74").Value = -63253.75
$ws.Range("H77").Value = 50041.332
$ws.Range("I77").Value = 64127.75
$ws.Range("K77").Value = 320638.75
$ws.Range("M77").Value = -316270.75
$ws.Range("H140").Value = 102806.336
$ws.Range("J140").Value = 102806.336
$ws.Range("L140").Value = 102806.336
$ws.Range("N140").Value = -113166.336

# --- CRP sheet updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1004.7
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 1180.875
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 1180.875
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = -1880.875
$ws.Range("H23").Value = 5000
$ws.Range("I23").Value = 5000
$ws.Range("K23").Value = 5000
$ws.Range("M23").Value = -4760
$ws.Range("H27").Value = 5000
$ws.Range("I27").Value = 5000
$ws.Range("K27").Value = 5000
$ws.Range("M27").Value = -4808
$ws.Range("H29").Value = 6000
$ws.Range("J29").Value = 6000
$ws.Range("L29").Value = 6000
$ws.Range("N29").Value = -6586
$ws.Range("H32").Value = 13727
$ws.Range("I32").Value = 13727
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 13727
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -13411
$ws.Range("N32").ClearContents()
$ws.Range("H33").Value = 25571.846
$ws.Range("I33").Value = 3746.5
$ws.Range("J33").Value = 44279.285
$ws.Range("K33").Value = 3746.5
$ws.Range("L33").Value = 44279.285
$ws.Range("M33").Value = -3367.5
$ws.Range("N33").Value = -45037.285
$ws.Range("H41").Value = 19999.8
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H105").Value = 1011.8947
$ws.Range("I105").Value = 1123.2307
$ws.Range("J105").Value = 770.6667
$ws.Range("K105").Value = 1123.2307
$ws.Range("L105").Value = 770.6667
$ws.Range("M105").Value = 623.7692999999999
$ws.Range("N105").Value = -4264.6667
$ws.Range("H134").Value = 11518.615
$ws.Range("I134").Value = 7092.1055
$ws.Range("J134").Value = 23533.428
$ws.Range("K134").Value = 21276.3165
$ws.Range("L134").Value = 70600.284
$ws.Range("M134").Value = -18741.3165
$ws.Range("N134").Value = -75670.284

# --- GSM sheet updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 279.83334
$ws.Range("I2").Value = 279.83334
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 279.83334
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -166.83334
$ws.Range("N2").ClearContents()
$ws.Range("H70").Value = 19315.666
$ws.Range("I70").Value = 18916.666
$ws.Range("J70").Value = 19448.666
$ws.Range("K70").Value = 18916.666
$ws.Range("L70").Value = 19448.666
$ws.Range("M70").Value = -18646.666
$ws.Range("N70").Value = -19988.666
$ws.Range("H73").Value = 19315.666
$ws.Range("I73").Value = 18916.666
$ws.Range("J73").Value = 19448.666
$ws.Range("K73").Value = 18916.666
$ws.Range("L73").Value = 19448.666
$ws.Range("M73").Value = -17980.666
$ws.Range("N73").Value = -21320.666

# --- LTW sheet updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 39500
$ws.Range("J94").Value = 39500
$ws.Range("L94").Value = 39500
$ws.Range("N94").Value = -40852
$ws.Range("H132").Value = 3099.64
$ws.Range("I132").Value = 2904.6191
$ws.Range("J132").Value = 4123.5
$ws.Range("K132").Value = 8713.8573
$ws.Range("L132").Value = 12370.5
$ws.Range("M132").Value = -6183.8573
$ws.Range("N132").Value = -17430.5
$ws.Range("H136").Value = 2959
$ws.Range("I136").Value = 2727.2144
$ws.Range("K136").Value = 8181.6432
$ws.Range("M136").Value = -5631.6432
$ws.Range("H140").Value = 134593.89
$ws.Range("J140").Value = 166916.17
$ws.Range("L140").Value = 166916.17
$ws.Range("N140").Value = -177276.17

# --- WVR sheet updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 14666.667
$ws.Range("I62").Value = 4000
$ws.Range("K62").Value = 4000
$ws.Range("M62").Value = -3376
$ws.Range("H65").Value = 14666.667
$ws.Range("I65").Value = 4000
$ws.Range("K65").Value = 20000
$ws.Range("M65").Value = -16880
$ws.Range("H81").Value = 7027.5
$ws.Range("I81").Value = 8812.799999999999
$ws.Range("J81").Value = 4593
$ws.Range("K81").Value = 17625.6
$ws.Range("L81").Value = 9186
$ws.Range("M81").Value = -16564.6
$ws.Range("N81").Value = -11308
$ws.Range("H84").Value = 7027.5
$ws.Range("I84").Value = 8812.799999999999
$ws.Range("J84").Value = 4593
$ws.Range("K84").Value = 88128
$ws.Range("L84").Value = 45930
$ws.Range("M84").Value = -82824
$ws.Range("N84").Value = -56538
